$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.493.12'
$ws.Range('E2').Value = '  -1.47%  '
$ws.Range('D3').Value = '2.056.66'
$ws.Range('E3').Value = '  +0.43%  '
$ws.Range('E4').Value = '  +0.11%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '242.25'
$c.ClearFormats()
$ws.Range('E5').Value = '  -1.78%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '0.662'
$c.ClearFormats()
$ws.Range('E6').Value = '  +0.52%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  -5.60%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '58.59'
$c.ClearFormats()
$ws.Range('E9').Value = '  -2.24%  '
$ws.Range('E10').Value = '  -5.55%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.0750'
$c.ClearFormats()
$ws.Range('E11').Value = '  -3.32%  '
$ws.Range('E12').Value = '  -3.04%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.905'
$c.ClearFormats()
$ws.Range('E13').Value = '  +1.23%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '14.70'
$c.ClearFormats()
$ws.Range('E14').Value = '  -6.54%  '
$ws.Range('E15').Value = '  +0.41%  '
$ws.Range('E16').Value = '  -5.92%  '
$ws.Range('D17').Value = '2.047.51'
$ws.Range('E17').Value = '  +0.50%  '
$ws.Range('D18').Value = '36.445.01'
$ws.Range('E18').Value = '  -1.42%  '
$ws.Range('E19').Value = '  -8.84%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '72.06'
$c.ClearFormats()
$ws.Range('E20').Value = '  -3.24%  '
$ws.Range('E21').Value = '  -4.95%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '237.67'
$c.ClearFormats()
$ws.Range('E22').Value = '  +0.64%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '5.23'
$c.ClearFormats()
$ws.Range('E23').Value = '  -4.26%  '
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('E25').Value = '  -4.54%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '9.31'
$c.ClearFormats()
$ws.Range('E26').Value = '  -2.43%  '
$ws.Range('E27').Value = '  -1.31%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '163.93'
$c.ClearFormats()
$ws.Range('E28').Value = '  -3.89%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '20.08'
$c.ClearFormats()
$ws.Range('E29').Value = '  +0.28%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '0.121'
$c.ClearFormats()
$ws.Range('E30').Value = '  -1.74%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '1.18'
$c.ClearFormats()
$ws.Range('E31').Value = '  +4.58%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '5.03'
$c.ClearFormats()
$ws.Range('E32').Value = '  -7.34%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '4.44'
$c.ClearFormats()
$ws.Range('E33').Value = '  -6.64%  '
$ws.Range('E34').Value = '  -4.01%  '
$ws.Range('E35').Value = '  +0.17%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '1.83'
$c.ClearFormats()
$ws.Range('E36').Value = '  -0.24%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '2.19'
$c.ClearFormats()
$ws.Range('E37').Value = '  -3.77%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.0821'
$c.ClearFormats()
$ws.Range('E38').Value = '  -5.69%  '
$ws.Range('E39').Value = '  -6.21%  '
$ws.Range('E40').Value = '  -4.31%  '
$ws.Range('E41').Value = '  -3.56%  '
$ws.Range('E42').Value = '  -8.82%  '
$ws.Range('E43').Value = '  -4.00%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '93.66'
$c.ClearFormats()
$ws.Range('E44').Value = '  -4.33%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.0904'
$c.ClearFormats()
$ws.Range('E45').Value = '  -8.66%  '
$ws.Range('D46').Value = '1.390.71'
$ws.Range('E46').Value = '  +7.11%  '
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '7.52'
$c.ClearFormats()
$ws.Range('E47').Value = '  +10.56%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '15.77'
$c.ClearFormats()
$ws.Range('E48').Value = '  -7.83%  '
$ws.Range('E49').Value = '  -0.46%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '2.26'
$c.ClearFormats()
$ws.Range('E50').Value = '  -4.90%  '
$ws.Range('D51').Value = '2.247.83'
$ws.Range('E51').Value = '  +0.75%  '
